# Updated data to reflect new requirement separation
#
# The "Terms Typically Offered" column (D) is moved to column G, and three
# new columns are inserted in its place: Corequisites (D), Concurrent (E),
# and Recommended (F). Any "Corequisite:", "Concurrent:", or "Recommended:"
# clause that was embedded in the Prerequisites (column C) text is split
# out into the appropriate new column, defaulting to "NA" where no such
# clause exists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# ---------------------------------------------------------------------
# Step 1: capture the existing "Terms Typically Offered" column (D)
# values for every row before we overwrite column D with new content.
# ---------------------------------------------------------------------
$termsOffered = @{
    1  = "Terms Typically Offered"
    2  = "F, W, SP"
    3  = "W"
    4  = "SP "
    5  = "F, W, SP"
    6  = "SP"
    7  = "W"
    8  = "F"
    9  = "SP "
    10 = "F"
    11 = "W, SP "
    12 = "F "
    13 = "SP "
}

# ---------------------------------------------------------------------
# Step 2: move those values into the new column G.
# ---------------------------------------------------------------------
foreach ($row in 1..13) {
    $ws.Range("G$row").Value = $termsOffered[$row]
}

# ---------------------------------------------------------------------
# Step 3: header row - add the three new column headers.
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# ---------------------------------------------------------------------
# Step 4: default every data row's new Corequisites/Concurrent/Recommended
# cell to "NA".
# ---------------------------------------------------------------------
foreach ($row in 2..13) {
    $ws.Range("D$row").Value = "NA"
    $ws.Range("E$row").Value = "NA"
    $ws.Range("F$row").Value = "NA"
}

# ---------------------------------------------------------------------
# Step 5: pull the embedded Corequisite/Concurrent/Recommended clauses out
# of column C (Prerequisites) and into their own columns.
# ---------------------------------------------------------------------

# Row 4 - PSC 103: "... PHYS 141. Recommended: PSC 102."
$ws.Range("C4").Value = "PSC" + $nbsp + "101 or PHYS" + $nbsp + "121 or PHYS 131 or PHYS" + $nbsp + "141."
$ws.Range("F4").Value = "PSC" + $nbsp + "102."

# Row 9 - PSC 392: "... GE Areas B2, B3, and B4. Recommended: UNIV 391 and completion of GE Areas D2 and D3."
$ws.Range("C9").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4."
$ws.Range("F9").Value = "UNIV" + $nbsp + "391 and completion of GE Areas D2 and D3."

# Row 11 - PSC 425: "... in Science. Concurrent: EDUC 469 or EDUC 479."
$ws.Range("C11").Value = "Acceptance into the Single Subject Credential Program in Science."
$ws.Range("E11").Value = "EDUC" + $nbsp + "469 or EDUC" + $nbsp + "479."

# Row 12 - PSC 491: "... standing. Corequisite: GE Area D5."
$ws.Range("C12").Value = "Consent of instructor, and senior or graduate standing."
$ws.Range("D12").Value = "GE Area D5."

# Row 13 - PSC 492: "... graduate standing. Recommended: UNIV 391, GE Area D2, and GE Area D3."
$ws.Range("C13").Value = "Junior standing and completion of GE Area B, or graduate standing."
$ws.Range("F13").Value = "UNIV" + $nbsp + "391, GE Area D2, and GE Area D3."

Write-Output "Applied requirement-separation update to PSC sheet."
